# Daily attendance processing - 2025-10-29 13:36:25
#
# Normalises the "Recorded By" column (G): whenever the literal "System"
# entry appears anywhere in the comma-separated list of recorders other
# than first, move it to the front while preserving the relative order of
# the remaining entries.
#
# Notes on this COM-interop runtime's quirks (discovered while iterating):
#   - Range/.Cells "Value" *getter* returns a stringified property
#     descriptor instead of the cell contents; use "Value2" (or "Text")
#     to read. "Value" as a *setter* works fine.
#   - "-ceq"/"-cmatch"/"-clike" do not reliably honour case sensitivity
#     here, so exact-case comparisons use System.String's .CompareTo().
#   - System.Collections.ArrayList's .ToArray() comes back empty, so
#     native PowerShell arrays (@() / +=) are used for rebuilding values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item(1048576, 1).End(-4162).Row
$recordedByCol = 7  # column G

$changed = 0

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $recordedByCol)
    $val = $cell.Value2

    if ($val -eq $null -or $val -eq "") {
        continue
    }

    $parts = @($val -split ", ")

    $sysCount = 0
    $sysIndex = -1
    for ($i = 0; $i -lt $parts.Count; $i++) {
        if ($parts[$i].CompareTo("System") -eq 0) {
            $sysCount++
            $sysIndex = $i
        }
    }

    # Only reorder when there is exactly one exact-case "System" token and
    # it is not already the first entry.
    if ($sysCount -eq 1 -and $sysIndex -ne 0) {
        $newParts = @("System")
        for ($i = 0; $i -lt $parts.Count; $i++) {
            if ($i -ne $sysIndex) {
                $newParts += $parts[$i]
            }
        }
        $newVal = $newParts -join ", "
        $cell.Value = $newVal
        $changed++
    }
}

Write-Host "Recorded By column normalised; rows changed:" $changed
